# Hjemme passive tweaks lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: subject/trial-count headers changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): D2 value removed, new values added at B2 and C2
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 26.835524999093739
$ws.Range("C2").Value = 25.532773574620265

# Row 3 (STR): B3 value removed, C3 value updated
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 23.893535323364024

# Selection narrowed to match the edited range
$ws.Range("B1:E3").Select()
